$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old numeric sample row (A1:C1 = 1,2,3) with a small
# "expense report" style table: header row + one data row, columns A:D.
$ws.Range("A1").Value = "人员"
$ws.Range("B1").Value = "费用"
$ws.Range("C1").Value = "经办人"
$ws.Range("D1").Value = "时间"

$ws.Range("A2").Value = "xx"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"

$ws.Range("C2").Value = "梅煜"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-10-09"

$ws.Range("D1").Select()
